{"js": "// Update the date line and the 25 \"two digit division\" answers in the table.\n// Each entry targets a specific (row, col) cell in the single table on the\n// page (row indices count every table row, including the blank spacer rows,\n// so only rows 0, 4, 8, 12, 16 are addressed). We match the exact original\n// text inside that cell and replace just that text range so the existing\n// run/paragraph formatting (fonts, size, alignment) is left untouched.\n\nconst tableEdits = [\n  [0, 0, \"99\u00f78=12, 3\", \"92\u00f78=11, 4\"],\n  [0, 1, \"27\u00f74=6, 3\", \"55\u00f72=27, 1\"],\n  [0, 2, \"74\u00f75=14, 4\", \"48\u00f75=9, 3\"],\n  [0, 3, \"74\u00f79=8, 2\", \"70\u00f74=17, 2\"],\n  [0, 4, \"15\u00f72=7, 1\", \"90\u00f74=22, 2\"],\n\n  [4, 0, \"65\u00f77=9, 2\", \"58\u00f73=19, 1\"],\n  [4, 1, \"56\u00f75=11, 1\", \"35\u00f78=4, 3\"],\n  [4, 2, \"72\u00f76=12, 0\", \"20\u00f77=2, 6\"],\n  [4, 3, \"81\u00f75=16, 1\", \"63\u00f73=21, 0\"],\n  [4, 4, \"59\u00f75=11, 4\", \"38\u00f74=9, 2\"],\n\n  [8, 0, \"48\u00f79=5, 3\", \"55\u00f73=18, 1\"],\n  [8, 1, \"16\u00f76=2, 4\", \"12\u00f72=6, 0\"],\n  [8, 2, \"58\u00f78=7, 2\", \"47\u00f72=23, 1\"],\n  [8, 3, \"90\u00f74=22, 2\", \"50\u00f72=25, 0\"],\n  [8, 4, \"67\u00f73=22, 1\", \"55\u00f78=6, 7\"],\n\n  [12, 0, \"94\u00f77=13, 3\", \"94\u00f76=15, 4\"],\n  [12, 1, \"22\u00f76=3, 4\", \"41\u00f75=8, 1\"],\n  [12, 2, \"43\u00f74=10, 3\", \"41\u00f74=10, 1\"],\n  [12, 3, \"75\u00f72=37, 1\", \"16\u00f77=2, 2\"],\n  [12, 4, \"89\u00f78=11, 1\", \"35\u00f72=17, 1\"],\n\n  [16, 0, \"66\u00f72=33, 0\", \"88\u00f77=12, 4\"],\n  [16, 1, \"83\u00f76=13, 5\", \"72\u00f74=18, 0\"],\n  [16, 2, \"94\u00f76=15, 4\", \"72\u00f72=36, 0\"],\n  [16, 3, \"63\u00f77=9, 0\", \"79\u00f72=39, 1\"],\n  [16, 4, \"93\u00f74=23, 1\", \"61\u00f77=8, 5\"],\n];\n\n// 1) Update the date heading above the table.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\nconst dateResults = datePara.search(\"2023-12-30 Saturday\", { matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2023-12-31 Sunday\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Update each answer cell in the table, one at a time, so a freshly\n// written value never gets confused with another cell's original text\n// (several answers coincide with other cells' old/new values).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of tableEdits) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Update the date line and the 25 \"two digit division\" answers in the table.\n# Cell/row indices below are 1-based (COM convention). Row numbers count\n# every table row, including the blank spacer rows, so only rows 1, 5, 9,\n# 13, 17 contain text. Assigning Range.Text in place keeps the existing\n# run/paragraph formatting (fonts, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading above the table.\n$d.Paragraphs.Item(1).Range.Text = \"2023-12-31 Sunday\"\n\n# 2) Table answers.\n$t = $d.Tables.Item(1)\n\n$tableEdits = @(\n    @(1, 1, \"92\u00f78=11, 4\"),\n    @(1, 2, \"55\u00f72=27, 1\"),\n    @(1, 3, \"48\u00f75=9, 3\"),\n    @(1, 4, \"70\u00f74=17, 2\"),\n    @(1, 5, \"90\u00f74=22, 2\"),\n\n    @(5, 1, \"58\u00f73=19, 1\"),\n    @(5, 2, \"35\u00f78=4, 3\"),\n    @(5, 3, \"20\u00f77=2, 6\"),\n    @(5, 4, \"63\u00f73=21, 0\"),\n    @(5, 5, \"38\u00f74=9, 2\"),\n\n    @(9, 1, \"55\u00f73=18, 1\"),\n    @(9, 2, \"12\u00f72=6, 0\"),\n    @(9, 3, \"47\u00f72=23, 1\"),\n    @(9, 4, \"50\u00f72=25, 0\"),\n    @(9, 5, \"55\u00f78=6, 7\"),\n\n    @(13, 1, \"94\u00f76=15, 4\"),\n    @(13, 2, \"41\u00f75=8, 1\"),\n    @(13, 3, \"41\u00f74=10, 1\"),\n    @(13, 4, \"16\u00f77=2, 2\"),\n    @(13, 5, \"35\u00f72=17, 1\"),\n\n    @(17, 1, \"88\u00f77=12, 4\"),\n    @(17, 2, \"72\u00f74=18, 0\"),\n    @(17, 3, \"72\u00f72=36, 0\"),\n    @(17, 4, \"79\u00f72=39, 1\"),\n    @(17, 5, \"61\u00f77=8, 5\")\n)\n\nforeach ($edit in $tableEdits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $newText = $edit[2]\n    $t.Cell($row, $col).Range.Text = $newText\n}\n"}
